$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Date value (B8) - updated publish date
$ws.Range("B8").Value = "2023-09-01T08:48:57+00:00"

# Case Sensitive value (B14) - was blank, now "true".
# Assigning the literal string "true" directly would be auto-coerced to a
# boolean by Value-type inference, so round-trip it through a formula that
# evaluates to the text "true", then paste the computed value back in place
# (keeping the existing cell style/format untouched).
$ws.Cells.Item(14, 2).Formula = '="true"'
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)
